$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("B18").Value = "ELM-2NA-Máquinas Elétricas"
$ws.Range("D18").Value = "-"
$ws.Range("E18").Value = "-"

# Row 19
$ws.Range("B19").Value = "ELM-2NA-Máquinas Elétricas"
$ws.Range("C19").Value = "[-, 'ELM-2NA-Instalções Elétricas']"
$ws.Range("E19").Value = "['ELM-2NA-Instalções Elétricas', -]"

# Row 20
$ws.Range("B20").Value = "-"
$ws.Range("C20").Value = "['ELM-2NA-Lab. De Máquinas elétricas', -]"
$ws.Range("E20").Value = "['ELM-2NA-Instalções Elétricas', -]"

# Row 21
$ws.Range("B21").Value = "-"
$ws.Range("C21").Value = "['ELM-2NA-Lab. De Máquinas elétricas', -]"
$ws.Range("E21").Value = "-"
